$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text (06/02/2023 ->
#    08/02/2023) on the Date Placeholder of the slide master and every
#    slide layout (Insert > Header & Footer > Apply to All equivalent).
# ---------------------------------------------------------------------------
$oldDate = "06/02/2023"
$newDate = "08/02/2023"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DateShape $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Rewrite the bullet text on slide 2's content placeholder.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

$tr2.Paragraphs(3, 1).Text = "The team consists of 5 members: Ben Kelly, Adrian Kucia, Kalina Filipowicz, Ed Davies and Charlie Callister."
$tr2.Paragraphs(4, 1).Text = "The requirements for this project were largely gathered via the initial meeting with the client, with further requirements, such as _____, coming to light as we further investigated what the client truly needed.  This investigation process allowed us to refine the identified requirements even further."

# Append two new, bullet-less blank paragraphs after the existing text.
$tr2.InsertAfter("`r`r")

$newPara1 = $tr2.Paragraphs(5, 1)
$newPara2 = $tr2.Paragraphs(6, 1)
$newPara1.ParagraphFormat.Bullet.Type = 0
$newPara2.ParagraphFormat.Bullet.Type = 0

Write-Output "done"
